# Applies the fix described in commit "Correction type pour génération à
# partir fsh ea4a6f04ed193a83290686b2f69a3f9cd2e7f4ad":
#  - Fills in the previously-empty "Name" value on the Metadata sheet with
#    "ProfessioncategoriesocioprofessionnelleVs"
#  - Refreshes the "Date" metadata value to reflect the regeneration time

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")

# Row 4 holds the "Name" property (A4="Name"); B4 was empty and now gets
# the value of the ValueSet name.
$metadata.Range("B4").Value = "ProfessioncategoriesocioprofessionnelleVs"

# Row 8 holds the "Date" property (A8="Date"); update its value (B8).
$metadata.Range("B8").Value = "2025-07-18T06:40:38+00:00"
